$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 13 new restaurant rows starting at row 8 (after existing 6 data rows + header)
$newRows = @(
    @('Yamasan Michishita Shoten','Hakodate','Hokkaido 040-0063','9-15 Wakamatsucho',41.7726921927474,140.72617457862),
    @('Masa Zushi','Otaru','Hokkaido 047-0024','1 Chome-1-1 Hanazono',44.680093291771797,140.890623885206),
    @('New Sanko Restaurant','Otaru','Hokkaido 047-0032','1-3-6 Inaho',43.195314008878398,140.99911878456399),
    @('Sapporo Ichiryuan','Sapporo','Hokkaido 060-0004 ',' B1, Hokuren Bldg., 1-1, Kitayonjo Nishi, Chuo-ku',43.067124839272097,141.354138867552),
    @('Soup Curry Shabazo','Sapporo','Hokkaido 060-0600 ','B1, Sapporo North Plaza, 4, Kitaichijo Nishi, Chuo-ku',43.116771913506902,141.35840197578199),
    @('Matsusakagyu Yakiniku','Chuo','Osaka 542-0076','1-1-19, Namba',34.692269665704401,135.50911722815599),
    @('Steakhouse Kozai','Chuo','Osaka 542-0083','1 Chome-17-15 Higashishinsaibashi',35.463249573212799,135.663244782601),
    @('Ali''s Kitchen Osaka Halal Restaurant','Chuo','Osaka 542-0085','1 Chome-10-12 Shinsaibashisuji',34.6748065646836,135.500938096878),
    @('Okonomiyaki Chitose','Nishinari Ward','Osaka 557-0002','1 Chome-11-10 Taishi',34.893970861434198,135.563846872672),
    @('Curry Yakumido','Haginochaya','Osaka 557-0004','2 Chome-2-10',34.909547149518403,135.436494507135),
    @('Himawari Tei','Fujikawaguchiko','Yamanashi 401-0301','7364-2 Funatsu',35.471925338743297,138.74964044721699),
    @('Sanrokuen','Fujikawaguchiko','Yamanashi 401-0301','3370-1 Funatsu',35.515319954042397,138.80705504684801),
    @('Hotokura Funari Fruitsline','Yamanashi','Yamanashi 405-0031','1091 Manriki',35.690321286489102,138.67492781028199),
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = [double]$row[4]
    $ws.Cells.Item($r, 6).Value = [double]$row[5]
}

# Sort the full data range (A2:F20) ascending by column C (Prefecture/zip)
$sortRange = $ws.Range("A2:F20")
$keyRange = $ws.Range("C2:C20")
$sortRange.Sort($keyRange)

# Update selection to match the post-edit saved state
$ws.Range("D23").Select()

